$d = $word.ActiveDocument

function Add-TrailingPeriod($paraIndex) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range.Duplicate
    $r.Collapse(0)        # wdCollapseEnd
    $r.MoveEnd(1, -1)      # wdCharacter: step back before the paragraph mark
    $r.InsertAfter(".")
}

# Paragraph indices (1-based, Word Paragraphs collection) that get a trailing "."
# appended right after their existing text.
$targets = @(2, 3, 5, 6, 8, 9, 10, 12, 13)
foreach ($idx in $targets) {
    Add-TrailingPeriod $idx
}

# Move the _GoBack bookmark from the start of the document to the end of the
# investment-loan paragraph (paragraph 13) -- i.e. the location of the last edit.
$d.Bookmarks("_GoBack").Delete()

$lastPara = $d.Paragraphs(13)

# Placing a zero-length bookmark exactly at "end of text, before the paragraph
# mark" needs a small trick: insert a temporary placeholder character after the
# text, anchor the bookmark just before the placeholder (a safe, non-boundary
# position), then remove the placeholder again. The bookmark stays put.
$tail = $lastPara.Range.Duplicate
$tail.Collapse(0)
$tail.MoveEnd(1, -1)
$tail.InsertAfter("X")

$afterText = $lastPara.Range.Duplicate
$afterText.SetRange($lastPara.Range.End - 2, $lastPara.Range.End - 2)
$d.Bookmarks.Add("_GoBack", $afterText)

$placeholder = $lastPara.Range.Duplicate
$placeholder.Collapse(0)
$placeholder.MoveEnd(1, -1)
$placeholder.MoveStart(1, -1)
$placeholder.Delete()
